$wb = $excel.ActiveWorkbook

# Sheets that contain the "想去人数" (interested count) column F that needs updating
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 10: F10 4892 -> 4893
    $ws.Range("F10").Value = 4893

    # Row 11: F11 4607 -> 4608
    $ws.Range("F11").Value = 4608
}
